$wb = $excel.ActiveWorkbook

# Update "想去人数" (column F) on sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 826
$ws1.Range("F3").Value = 4282
$ws1.Range("F4").Value = 120
$ws1.Range("F5").Value = 770

# Update "想去人数" (column F) on sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 826
$ws4.Range("F3").Value = 4282
$ws4.Range("F4").Value = 120
$ws4.Range("F5").Value = 770
